# Ispravak pravopisa i gramatike
# Applies a series of targeted Find & Replace corrections across the document.

$d = $word.ActiveDocument

function Replace-Text {
    param(
        [string]$Find,
        [string]$Replace
    )
    $rng = $d.Content
    $ok = $rng.Find.Execute($Find, $true, $false, $false, $false, $false, $true, 1, $false, $Replace, 2)
    if (-not $ok) {
        Write-Host "WARNING: not found -> $Find"
    }
}

Replace-Text "S obzirom da je naš zadatak" "S obzirom na to da je naš zadatak"
Replace-Text "razgovara sa modelom" "razgovara s modelom"
Replace-Text "priča izravno sa bazom podataka ili baza" "priča izravno s bazom podataka ili baza"
Replace-Text "datumRegistracije: DATE" "datumRegistracije: TIMESTAMP"
Replace-Text "jeRegistriran: BOOLEAN" "jeValidiran: BOOLEAN"
Replace-Text "brojTelefona. U bazi podata se ne pamti " "brojTelefona, je li korisnik verificirao svoj račun. U bazi podata se ne pamti "
Replace-Text "datum kada je posalana" "datum kada je poslana"
Replace-Text "ukNovacPotrosen: INT - ukupna kolicina novaca koju je korisnik potrosio" "ukNovacPotrosen: INT - ukupna količina novaca koju je korisnik potrošio"
Replace-Text "ukNovacZaraden: INT - ukupna kolicina novaca koju je korisnik zaradio" "ukNovacZaraden: INT - ukupna količina novaca koju je korisnik zaradio"
Replace-Text "prosjecnaOcjena: Double - prosjecna ocjena korisnika" "prosjecnaOcjena: Double - prosječna ocjena korisnika"
Replace-Text " - vrijeme kada ce se posao održ" " - vrijeme kada će se posao održ"
Replace-Text "trajanje: SMALLINT - okvirno vrijeme potrebno za odradivanje posla" "trajanje: LONG - okvirno vrijeme potrebno za odrađivanje posla"
Replace-Text "ponudeniNovac: SMALLINT - novac ponuden za odradivanje posla" "ponudeniNovac: INT - novac ponuđen za odrađivanje posla"
Replace-Text "posaoGotov: BOOLEAN - zastavica je li posao odraden" "posaoGotov: BOOLEAN - zastavica je li posao odrađen"
Replace-Text "255) - opširan opis što sve kategorija podrazumjeva" "255) - opširan opis što sve kategorija podrazumijeva"
Replace-Text " - treuntak kada je posao odrađen" " - trenutak kada je posao odrađen"
Replace-Text "korisnike (posloprimce i poslodavce sa poslom)" "korisnike (posloprimce i poslodavce s poslom)"

Write-Host "Done"
